# Apply cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "68.266.48", "  +1.55%  ", 0)
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.925.98", "  -0.44%  ", 0)
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.03%  ", 1)
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "487.30", "  +3.34%  ", 1)
    ,@(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "148.32", "  +1.40%  ", 1)
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.624", "  -0.22%  ", 1)
    ,@(8, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.998", "  -0.09%  ", 1)
    ,@(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.734", "  +0.02%  ", 1)
    ,@(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.167", "  +1.64%  ", 1)
    ,@(11, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000349", "  +4.19%  ", 1)
    ,@(12, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "43.13", "  -0.93%  ", 1)
    ,@(13, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "10.74", "  +3.37%  ", 1)
    ,@(14, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "4.548.62", "  -0.55%  ", 0)
    ,@(15, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.921.05", "  -1.07%  ", 0)
    ,@(16, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "14.45", "  -5.11%  ", 1)
    ,@(17, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.136", "  -0.75%  ", 1)
    ,@(18, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "19.95", "  +0.43%  ", 1)
    ,@(19, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.14", "  -2.10%  ", 1)
    ,@(20, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "68.355.41", "  +1.18%  ", 0)
    ,@(21, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "442.58", "  +0.98%  ", 1)
    ,@(22, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "15.17", "  +4.50%  ", 1)
    ,@(23, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "3.48", "  +2.47%  ", 1)
    ,@(24, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "88.21", "  +0.68%  ", 1)
    ,@(25, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "11.32", "  +15.15%  ", 1)
    ,@(26, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "11.45", "  +11.41%  ", 1)
    ,@(27, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.63", "  +0.63%  ", 1)
    ,@(28, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "38.74", "  -0.60%  ", 1)
    ,@(29, "LEO", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo", "5.71", "  -0.94%  ", 1)
    ,@(30, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "718.31", "  -0.73%  ", 1)
    ,@(31, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "13.77", "  +1.35%  ", 1)
    ,@(32, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.131", "  -1.01%  ", 1)
    ,@(33, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.91", "  +3.26%  ", 1)
    ,@(34, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "6.45", "  +19.83%  ", 1)
    ,@(35, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "42.55", "  -0.89%  ", 1)
    ,@(36, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0887", "  +13.65%  ", 0)
    ,@(37, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "60.87", "  +5.21%  ", 1)
    ,@(38, "TheGraph", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt", "0.413", "  +22.33%  ", 1)
    ,@(39, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.150", "  -1.98%  ", 1)
    ,@(40, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "3.04", "  +16.94%  ", 1)
    ,@(41, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.07%  ", 1)
    ,@(42, "ThetaToken", "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta", "3.28", "  +6.96%  ", 1)
    ,@(43, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0482", "  +0.84%  ", 1)
    ,@(44, "WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "2.91", "  +3.38%  ", 1)
    ,@(45, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.142", "  +0.51%  ", 1)
    ,@(46, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "1.00", "  +0.01%  ", 1)
    ,@(47, "ApeXProtocol", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex", "3.32", "  +5.22%  ", 1)
    ,@(48, "BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.0₆0358", "  +36.84%  ", 0)
    ,@(49, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "3.43", "  -1.35%  ", 1)
    ,@(50, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "2.14", "  -2.04%  ", 1)
    ,@(51, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "145.99", "  -0.89%  ", 1)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $dCell = $ws.Cells.Item($r, 4)
    if ($item[5] -eq 1) {
        $dCell.NumberFormat = "@"
    }
    $dCell.Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
}
